$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the report title strings (new volume/issue number + new week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  9"
$ws.Range("C9").Value = "Report Covering the Week  2/24/2025  Through  3/2/2025"

# --- Update weekly crime statistics table (rows 15-31) ---

# Row 15
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0"
$ws.Range("C29").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "***.*"
$ws.Range("E29").Copy()
$ws.Range("H15").PasteSpecial(-4122)

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -57.894736842105
$ws.Range("I16").Value = 19
$ws.Range("J16").Value = 30
$ws.Range("K16").Value = -36.666666666666
$ws.Range("L16").Value = -9.523809523809
$ws.Range("M16").Value = -44.117647058823
$ws.Range("N16").Value = -84.297520661157

# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 40
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 90
$ws.Range("I17").Value = 35
$ws.Range("J17").Value = 23
$ws.Range("K17").Value = 52.173913043478
$ws.Range("L17").Value = 59.090909090909
$ws.Range("M17").Value = 75
$ws.Range("N17").Value = -42.622950819672

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 24
$ws.Range("J18").Value = 26
$ws.Range("K18").Value = -7.692307692307
$ws.Range("L18").Value = 20
$ws.Range("M18").Value = -4
$ws.Range("N18").Value = -87.368421052631

# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -42.857142857142
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = -29.729729729729
$ws.Range("I19").Value = 64
$ws.Range("J19").Value = 72
$ws.Range("K19").Value = -11.111111111111
$ws.Range("L19").Value = -29.670329670329
$ws.Range("M19").Value = -17.948717948717
$ws.Range("N19").Value = -59.748427672956

# Row 20
$ws.Range("I29").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 1
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("C29").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$ws.Range("E29").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("I29").Copy()
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -75
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 8
$ws.Range("K20").Value = -50
$ws.Range("L20").Value = -80
$ws.Range("M20").Value = -20
$ws.Range("N20").Value = -97.402597402597

# Row 21
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -8.333333333333
$ws.Range("F21").Value = 64
$ws.Range("G21").Value = 82
$ws.Range("H21").Value = -21.951219512195
$ws.Range("I21").Value = 151
$ws.Range("J21").Value = 162
$ws.Range("K21").Value = -6.79012345679
$ws.Range("L21").Value = -15.642458100558
$ws.Range("M21").Value = -7.926829268292
$ws.Range("N21").Value = -78.210678210678

# Row 22
$ws.Range("I29").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 2
$ws.Range("L29").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = -50
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = 8
$ws.Range("K22").Value = -37.5
$ws.Range("L22").Value = 66.666666666666
$ws.Range("M22").Value = 66.666666666666

# Row 23
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 133.333333333333
$ws.Range("F23").Value = 16
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 30
$ws.Range("J23").Value = 16
$ws.Range("K23").Value = 87.5
$ws.Range("L23").Value = 66.666666666666
$ws.Range("M23").Value = 57.894736842105

# Row 24
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 28
$ws.Range("F24").Value = 169
$ws.Range("G24").Value = 102
$ws.Range("H24").Value = 65.686274509803
$ws.Range("I24").Value = 296
$ws.Range("J24").Value = 226
$ws.Range("K24").Value = 30.973451327433
$ws.Range("L24").Value = -13.953488372093
$ws.Range("M24").Value = 88.535031847133

# Row 25
$ws.Range("C25").Value = 21
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 117
$ws.Range("G25").Value = 73
$ws.Range("H25").Value = 60.273972602739
$ws.Range("I25").Value = 199
$ws.Range("J25").Value = 150
$ws.Range("K25").Value = 32.666666666666
$ws.Range("L25").Value = -15.677966101694

# Row 26
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 100
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = 87.5
$ws.Range("I26").Value = 55
$ws.Range("J26").Value = 45
$ws.Range("K26").Value = 22.222222222222
$ws.Range("L26").Value = 27.906976744186
$ws.Range("M26").Value = 12.244897959183

# Row 27
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "0"
$ws.Range("C29").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "***.*"
$ws.Range("E29").Copy()
$ws.Range("H27").PasteSpecial(-4122)

# Row 28
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C29").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("I29").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 1
$ws.Range("L29").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 66.666666666666
$ws.Range("J28").Value = 5
$ws.Range("K28").Value = 80

# Row 31
$ws.Range("F31").Value = 1
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "0"
$ws.Range("C29").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "***.*"
$ws.Range("E29").Copy()
$ws.Range("H31").PasteSpecial(-4122)
